# Removed unused containers from azure blob storage map excel spreadsheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Containers")

# Delete the "sessionsummaries" and "sessionanalyses" rows (original rows 3-4),
# and the "useraccounts" row (original row 5) -- all unused/leftover containers --
# shifting the remaining rows (userphotos, activitylogs, messagesubmissions) up.
$ws.Rows("3:5").Delete() | Out-Null

# Reset the selection back to the top of the sheet now that the old selection
# (row 7) no longer points at meaningful data.
$ws.Range("A1").Select() | Out-Null
